# March 24 update 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns (same header style as the existing header row)
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"
$ws.Range("M1:O1").Font.Bold = $true
$ws.Range("M1:O1").HorizontalAlignment = -4108
$ws.Range("M1:O1").VerticalAlignment = -4160
$ws.Range("M1:O1").Borders.LineStyle = 1

# Fill new columns for each data row (2..26)
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"
    $ws.Cells.Item($r, 14).Value = 20140231
    $ws.Cells.Item($r, 15).Value = 5
}
